$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated NATMI Wnt5a-Fzd1 results (per Dr Hou's advice): recomputed stats
# for the existing ECs->{ECs,FAPs,M1,M2,sCs} rows (2-6) and new rows (7-11)
# for FAPs->{ECs,FAPs,M1,M2,sCs}.
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target
# cluster, E..T the numeric NATMI metrics.
$data = @(
  @("ECs", "Wnt5a", "Fzd1", "ECs", 1, 0.3333333333333333, 0.1599003333333333, 0.479701, 0.0264777194346773, 0.02647771943467731, 3, 1, 1.269979666666667, 3.809939, 0.05093337778110388, 0.05093337778110388, 0.2030701720265556, 1.827631548239, 0.001348599686748495, 0.001348599686748496),
  @("ECs", "Wnt5a", "Fzd1", "FAPs", 1, 0.3333333333333333, 0.1599003333333333, 0.479701, 0.0264777194346773, 0.02647771943467731, 3, 1, 20.604156, 61.812468, 0.8263433572627787, 0.8263433572627789, 3.294611412452, 29.651502712068, 0.02187968757031317, 0.02187968757031317),
  @("ECs", "Wnt5a", "Fzd1", "M1", 1, 0.3333333333333333, 0.1599003333333333, 0.479701, 0.0264777194346773, 0.02647771943467731, 3, 1, 0.6519523333333334, 1.955857, 0.02614698121592406, 0.02614698121592406, 0.1042473954174444, 0.938226558757, 0.0006923124326990148, 0.000692312432699015),
  @("ECs", "Wnt5a", "Fzd1", "M2", 1, 0.3333333333333333, 0.1599003333333333, 0.479701, 0.0264777194346773, 0.02647771943467731, 3, 1, 0.7481240000000001, 2.244372, 0.03000400976428538, 0.03000400976428539, 0.1196252769746667, 1.076627492772, 0.0007944377524540667, 0.0007944377524540669),
  @("ECs", "Wnt5a", "Fzd1", "sCs", 1, 0.3333333333333333, 0.1599003333333333, 0.479701, 0.0264777194346773, 0.02647771943467731, 3, 1, 1.659922, 4.979766, 0.06657227397590788, 0.0665722739759079, 0.2654220811073333, 2.388798729966, 0.001762681992462558, 0.001762681992462559),
  @("FAPs", "Wnt5a", "Fzd1", "ECs", 3, 1, 5.879152, 17.637456, 0.9735222805653226, 0.9735222805653228, 3, 1, 1.269979666666667, 3.809939, 0.05093337778110388, 0.05093337778110388, 7.466403497242667, 67.197631475184, 0.04958477809435538, 0.04958477809435539),
  @("FAPs", "Wnt5a", "Fzd1", "FAPs", 3, 1, 5.879152, 17.637456, 0.9735222805653226, 0.9735222805653228, 3, 1, 20.604156, 61.812468, 0.8263433572627787, 0.8263433572627789, 121.134964955712, 1090.214684601408, 0.8044636696924655, 0.8044636696924657),
  @("FAPs", "Wnt5a", "Fzd1", "M1", 3, 1, 5.879152, 17.637456, 0.9735222805653226, 0.9735222805653228, 3, 1, 0.6519523333333334, 1.955857, 0.02614698121592406, 0.02614698121592406, 3.832926864421334, 34.496341779792, 0.02545466878322504, 0.02545466878322505),
  @("FAPs", "Wnt5a", "Fzd1", "M2", 3, 1, 5.879152, 17.637456, 0.9735222805653226, 0.9735222805653228, 3, 1, 0.7481240000000001, 2.244372, 0.03000400976428538, 0.03000400976428539, 4.398334710848001, 39.58501239763201, 0.02920957201183131, 0.02920957201183132),
  @("FAPs", "Wnt5a", "Fzd1", "sCs", 3, 1, 5.879152, 17.637456, 0.9735222805653226, 0.9735222805653228, 3, 1, 1.659922, 4.979766, 0.06657227397590788, 0.0665722739759079, 9.758933746144, 87.830403715296, 0.06480959198344532, 0.06480959198344534)
)

$rowStart = 2
for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $rowStart + $i
  for ($j = 0; $j -lt 20; $j++) {
    $col = $j + 1
    $ws.Cells.Item($r, $col).Value = $data[$i][$j]
  }
}

